# Fixed some inconsistencies between the scrum board and the burndown chart
# (Sprint3 sheet): correct a few logged-work entries so the "DONE" /
# remaining-work formulas line up with the real numbers, which also
# cascades into the daily-total / remaining-work summary rows used by
# the burndown chart series.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Burndown Chart Sprint3")

# Task row 8: K8 logged work 0.15 -> 0.1
$ws.Range("K8").Value = 0.1

# Task row 13: F13 and G13 logged work 0.25 -> 0.2
$ws.Range("F13").Value = 0.2
$ws.Range("G13").Value = 0.2

# L8/N8, L13/N13 and the daily-total / remaining-work rows (38-39) are all
# formula-driven off the cells above, so Excel recalculates them
# automatically once the inputs change.
